# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310
#   *_new -> *_FV2404
# Columns A-J (1-10): "<name>_old"  -> "<name>_FV2310"
# Column K (11): "diff" stays unchanged
# Columns L-U (12-21): "<name>_new" -> "<name>_FV2404"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2310 = "_FV2310"
$fv2404 = "_FV2404"

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = $cell.Value2
    if ($header -like "*$oldSuffix") {
        $base = $header.Substring(0, $header.Length - $oldSuffix.Length)
        $cell.Value = $base + $fv2310
    } elseif ($header -like "*$newSuffix") {
        $base = $header.Substring(0, $header.Length - $newSuffix.Length)
        $cell.Value = $base + $fv2404
    }
}

# Freeze the header row (split below row 1) - matches the new <pane>/<selection>
# entries added to the sheetView in the target workbook.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the header range + data into an Excel Table ("Table1") so the headers
# get the table's autofilter / tableParts wiring the diff adds.
$usedRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"
